$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-24 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-25 Tuesday", 2)

$d.Content.Find.Execute("399÷5=79, 4", $true, $false, $false, $false, $false, $true, 1, $false, "753÷5=150, 3", 2)
$d.Content.Find.Execute("654÷6=109, 0", $true, $false, $false, $false, $false, $true, 1, $false, "862÷9=95, 7", 2)
$d.Content.Find.Execute("215÷4=53, 3", $true, $false, $false, $false, $false, $true, 1, $false, "692÷9=76, 8", 2)
$d.Content.Find.Execute("219÷2=109, 1", $true, $false, $false, $false, $false, $true, 1, $false, "801÷2=400, 1", 2)
$d.Content.Find.Execute("256÷5=51, 1", $true, $false, $false, $false, $false, $true, 1, $false, "544÷2=272, 0", 2)

$d.Content.Find.Execute("497÷8=62, 1", $true, $false, $false, $false, $false, $true, 1, $false, "530÷4=132, 2", 2)
$d.Content.Find.Execute("876÷9=97, 3", $true, $false, $false, $false, $false, $true, 1, $false, "477÷5=95, 2", 2)
$d.Content.Find.Execute("777÷6=129, 3", $true, $false, $false, $false, $false, $true, 1, $false, "435÷7=62, 1", 2)
$d.Content.Find.Execute("385÷4=96, 1", $true, $false, $false, $false, $false, $true, 1, $false, "318÷9=35, 3", 2)
$d.Content.Find.Execute("561÷2=280, 1", $true, $false, $false, $false, $false, $true, 1, $false, "951÷6=158, 3", 2)

$d.Content.Find.Execute("348÷4=87, 0", $true, $false, $false, $false, $false, $true, 1, $false, "554÷6=92, 2", 2)
$d.Content.Find.Execute("342÷2=171, 0", $true, $false, $false, $false, $false, $true, 1, $false, "800÷4=200, 0", 2)
$d.Content.Find.Execute("774÷3=258, 0", $true, $false, $false, $false, $false, $true, 1, $false, "617÷8=77, 1", 2)
$d.Content.Find.Execute("511÷9=56, 7", $true, $false, $false, $false, $false, $true, 1, $false, "244÷8=30, 4", 2)
$d.Content.Find.Execute("552÷8=69, 0", $true, $false, $false, $false, $false, $true, 1, $false, "969÷8=121, 1", 2)

$d.Content.Find.Execute("464÷2=232, 0", $true, $false, $false, $false, $false, $true, 1, $false, "714÷7=102, 0", 2)
$d.Content.Find.Execute("624÷8=78, 0", $true, $false, $false, $false, $false, $true, 1, $false, "546÷7=78, 0", 2)
$d.Content.Find.Execute("564÷4=141, 0", $true, $false, $false, $false, $false, $true, 1, $false, "794÷3=264, 2", 2)
$d.Content.Find.Execute("372÷7=53, 1", $true, $false, $false, $false, $false, $true, 1, $false, "866÷9=96, 2", 2)
$d.Content.Find.Execute("972÷9=108, 0", $true, $false, $false, $false, $false, $true, 1, $false, "188÷4=47, 0", 2)

$d.Content.Find.Execute("679÷7=97, 0", $true, $false, $false, $false, $false, $true, 1, $false, "984÷4=246, 0", 2)
$d.Content.Find.Execute("362÷4=90, 2", $true, $false, $false, $false, $false, $true, 1, $false, "749÷7=107, 0", 2)
$d.Content.Find.Execute("924÷5=184, 4", $true, $false, $false, $false, $false, $true, 1, $false, "853÷3=284, 1", 2)
$d.Content.Find.Execute("732÷2=366, 0", $true, $false, $false, $false, $false, $true, 1, $false, "471÷8=58, 7", 2)
$d.Content.Find.Execute("410÷4=102, 2", $true, $false, $false, $false, $false, $true, 1, $false, "478÷2=239, 0", 2)
